# Apply the update to the training schedule sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the row-2 values that changed for this condition.
$ws.Range("D2").Value = 4
$ws.Range("F2").Value = -3
$ws.Range("H2").Value = 46

# Move the active selection to E2 (was D5).
$ws.Range("E2").Select()
